# Apply the edit described by the diff:
#  1. Remove the "Meta description: ..." paragraph that follows the H1 title.
#  2. Replace the trailing italic "image prompt" paragraph with two new
#     paragraphs: a bold title line and an italic meta-description line.

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (paragraph 2) ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- Step 2: replace the final paragraph (the AI image-prompt text) ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newXml = "<w:p $ns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Attila Slot for Free - Review &amp; Ratings 2021</w:t></w:r></w:p>" + `
          "<w:p $ns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Attila online slot game and play for free. Discover its bonuses and immersive graphics. Ratings &amp; winning potential info.</w:t></w:r></w:p>"

$lastPara.Range.InsertXML($newXml) | Out-Null
